# Natmi following Dr Hou advice
# Update LR-pair summary stats (Agrn-Musk) to the full 4x4 sending/target cluster
# cross-tabulation (ECs / FAPs / M2 / sCs), replacing the 8-row partial table with
# the revised 16-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowData = @("ECs", "Agrn", "Musk", "ECs", 3, 1, 18.45425133333334, 55.362754, 0.3394975357727733, 0.3394975357727733, 2, 0.6666666666666666, 0.5981746666666666, 1.794524, 0.02144965888565845, 0.02144965888565845, 11.03886563989956, 99.34979075909601, 0.007282106334847611, 0.007282106334847613)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(2, $col).Value = $rowData[$col - 1]
}

$rowData = @("ECs", "Agrn", "Musk", "FAPs", 3, 1, 18.45425133333334, 55.362754, 0.3394975357727733, 0.3394975357727733, 3, 1, 6.923267333333333, 20.769802, 0.2482581275160803, 0.2482581275160804, 127.7637154171898, 1149.873438754708, 0.08428302252727218, 0.0842830225272722)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(3, $col).Value = $rowData[$col - 1]
}

$rowData = @("ECs", "Agrn", "Musk", "M2", 3, 1, 18.45425133333334, 55.362754, 0.3394975357727733, 0.3394975357727733, 2, 0.6666666666666666, 0.1586656666666667, 0.475997, 0.005689516150576289, 0.00568951615057629, 2.928056090637556, 26.352504815738, 0.001931576712860045, 0.001931576712860045)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(4, $col).Value = $rowData[$col - 1]
}

$rowData = @("ECs", "Agrn", "Musk", "sCs", 3, 1, 18.45425133333334, 55.362754, 0.3394975357727733, 0.3394975357727733, 3, 1, 20.20726666666667, 60.6218, 0.7246026974476848, 0.7246026974476849, 372.9099778263556, 3356.1898004372, 0.2460008301977934, 0.2460008301977934)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(5, $col).Value = $rowData[$col - 1]
}

$rowData = @("FAPs", "Agrn", "Musk", "ECs", 3, 1, 15.88630666666666, 47.65891999999999, 0.2922557988641919, 0.292255798864192, 2, 0.6666666666666666, 0.5981746666666666, 1.794524, 0.02144965888565845, 0.02144965888565845, 9.502786194897777, 85.52507575407999, 0.006268787192992521, 0.006268787192992523)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(6, $col).Value = $rowData[$col - 1]
}

$rowData = @("FAPs", "Agrn", "Musk", "FAPs", 3, 1, 15.88630666666666, 47.65891999999999, 0.2922557988641919, 0.292255798864192, 3, 1, 6.923267333333333, 20.769802, 0.2482581275160803, 0.2482581275160804, 109.9851479926489, 989.8663319338398, 0.07255487738174048, 0.0725548773817405)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(7, $col).Value = $rowData[$col - 1]
}

$rowData = @("FAPs", "Agrn", "Musk", "M2", 3, 1, 15.88630666666666, 47.65891999999999, 0.2922557988641919, 0.292255798864192, 2, 0.6666666666666666, 0.1586656666666667, 0.475997, 0.005689516150576289, 0.00568951615057629, 2.520611438137778, 22.68550294324, 0.001662794087737395, 0.001662794087737396)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(8, $col).Value = $rowData[$col - 1]
}

$rowData = @("FAPs", "Agrn", "Musk", "sCs", 3, 1, 15.88630666666666, 47.65891999999999, 0.2922557988641919, 0.292255798864192, 3, 1, 20.20726666666667, 60.6218, 0.7246026974476848, 0.7246026974476849, 321.0188351617777, 2889.169516456, 0.2117693402017215, 0.2117693402017215)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(9, $col).Value = $rowData[$col - 1]
}

$rowData = @("M2", "Agrn", "Musk", "ECs", 3, 1, 4.826026, 14.478078, 0.08878300750222795, 0.08878300750222796, 2, 0.6666666666666666, 0.5981746666666666, 1.794524, 0.02144965888565845, 0.02144965888565845, 2.886806493874666, 25.981258444872, 0.001904365225765644, 0.001904365225765645)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(10, $col).Value = $rowData[$col - 1]
}

$rowData = @("M2", "Agrn", "Musk", "FAPs", 3, 1, 4.826026, 14.478078, 0.08878300750222795, 0.08878300750222796, 3, 1, 6.923267333333333, 20.769802, 0.2482581275160803, 0.2482581275160804, 33.41186815561733, 300.706813400556, 0.02204110319774922, 0.02204110319774923)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(11, $col).Value = $rowData[$col - 1]
}

$rowData = @("M2", "Agrn", "Musk", "M2", 3, 1, 4.826026, 14.478078, 0.08878300750222795, 0.08878300750222796, 2, 0.6666666666666666, 0.1586656666666667, 0.475997, 0.005689516150576289, 0.00568951615057629, 0.7657246326406667, 6.891521693766, 0.0005051323550806617, 0.0005051323550806619)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(12, $col).Value = $rowData[$col - 1]
}

$rowData = @("M2", "Agrn", "Musk", "sCs", 3, 1, 4.826026, 14.478078, 0.08878300750222795, 0.08878300750222796, 3, 1, 20.20726666666667, 60.6218, 0.7246026974476848, 0.7246026974476849, 97.52079432226665, 877.6871489004, 0.06433240672363241, 0.06433240672363243)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(13, $col).Value = $rowData[$col - 1]
}

$rowData = @("sCs", "Agrn", "Musk", "ECs", 3, 1, 15.19095733333334, 45.572872, 0.2794636578608068, 0.2794636578608069, 2, 0.6666666666666666, 0.5981746666666666, 1.794524, 0.02144965888565845, 0.02144965888565845, 9.086845839214222, 81.78161255292801, 0.005994400132052667, 0.00599440013205267)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(14, $col).Value = $rowData[$col - 1]
}

$rowData = @("sCs", "Agrn", "Musk", "FAPs", 3, 1, 15.19095733333334, 45.572872, 0.2794636578608068, 0.2794636578608069, 3, 1, 6.923267333333333, 20.769802, 0.2482581275160803, 0.2482581275160804, 105.1710586679271, 946.539528011344, 0.06937912440931844, 0.06937912440931845)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(15, $col).Value = $rowData[$col - 1]
}

$rowData = @("sCs", "Agrn", "Musk", "M2", 3, 1, 15.19095733333334, 45.572872, 0.2794636578608068, 0.2794636578608069, 2, 0.6666666666666666, 0.1586656666666667, 0.475997, 0.005689516150576289, 0.00568951615057629, 2.410283372598223, 21.692550353384, 0.001590012994898187, 0.001590012994898187)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(16, $col).Value = $rowData[$col - 1]
}

$rowData = @("sCs", "Agrn", "Musk", "sCs", 3, 1, 15.19095733333334, 45.572872, 0.2794636578608068, 0.2794636578608069, 3, 1, 20.20726666666667, 60.6218, 0.7246026974476848, 0.7246026974476849, 306.9677257566223, 2762.7095318096, 0.2025001203245375, 0.2025001203245376)
for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item(17, $col).Value = $rowData[$col - 1]
}
